# Issue #13: allow two columns in metadata files to be related to create
# hierarchical SKOS. A new row is inserted right below the header row
# (new row 2) holding short "slug" identifiers for every column; the rows
# that used to be 2-5 (full iaest-measure/iaest-dimension ids, medida/dim,
# xsd:int/skos:Concept/URI-Comunidad, mapping-*.xlsx) all shift down by one
# row to become rows 3-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2 - pushes the old rows 2..5 down to 3..6
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the short slug identifiers
$ws.Range("A2").Value = "numero-de-edificios"
$ws.Range("B2").Value = "evacuacion-aguas-residuales"
$ws.Range("C2").Value = "agua-corriente"
$ws.Range("D2").Value = "portero-automatico"
$ws.Range("E2").Value = "agua-caliente-central"
$ws.Range("F2").Value = "ascensor"
$ws.Range("G2").Value = "portero"
$ws.Range("H2").Value = "gas"
$ws.Range("I2").Value = "aragon"
$ws.Range("J2").Value = "accesible"
$ws.Range("K2").Value = "telefono"
